# Update the "Förändrad" date column (C) for rows 2-20 from 45185 (2023-09-16)
# to 45204 (2023-10-05), keeping the existing date formatting/style untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 20; $row++) {
    $ws.Cells.Item($row, 3).Value = 45204
}
